$d = $word.ActiveDocument

# --- 1) Insert new blank heading-style paragraph after the section-break paragraph,
#        right before "HOP DONG LAO DONG" ---
$secBreakPara = $d.Paragraphs(6)
$secBreakPara.Range.InsertParagraphAfter()
$newBlank = $d.Paragraphs(7)
$newBlank.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="60"/><w:ind w:left="720" w:hanging="720"/><w:jc w:val="center"/><w:rPr><w:b/><w:color w:val="000000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p>')

# --- 2) "- Thoi gian lam viec: {JobType}" paragraph: drop trailing green space-run,
#        recolor paragraph-mark to black ---
$pThoiGian = $d.Paragraphs(25)
$pThoiGian.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="60"/><w:ind w:firstLine="709"/><w:jc w:val="both"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>- Thời</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> gian</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> làm việc:</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> {JobType}</w:t></w:r></w:p>')

# --- 3) "- Dia diem lam viec: {Place}" -> "- Dia diem lam viec: {Type}" ---
$pDiaDiem = $d.Paragraphs(26)
$pDiaDiem.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="60"/><w:ind w:firstLine="709"/><w:jc w:val="both"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Địa điểm làm việc</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>: {Type}</w:t></w:r></w:p>')

# --- 4) "- Cong viec phai lam: {Jobs}" -> "- Cong viec phai lam: {C}" (bookmark removed) ---
$pCongViec = $d.Paragraphs(27)
$pCongViec.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="60"/><w:ind w:firstLine="709"/><w:jc w:val="both"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">- Công việc phải làm: </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>{C}</w:t></w:r></w:p>')

# --- 5) "Dieu 3." heading paragraph: spacing 120/120 -> 60/60, add firstLine ind 709 ---
$pDieu3 = $d.Paragraphs(28)
$pDieu3.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="60"/><w:ind w:firstLine="709"/><w:jc w:val="both"/><w:rPr><w:b/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="000000"/></w:rPr><w:t>Điều 3. Các quyền lợi, quyền hạn và nghĩa vụ của Bên B:</w:t></w:r></w:p>')

Write-Host "edit complete"
